# [LPF-879]: CCMS Third party report
#
# The template used to ship with three tabs: MAIN (the report itself) plus
# two "source data" worksheets, DATA and "Adjusted Expenditure", that MAIN's
# formulas pulled from. Those two tabs are no longer needed in the template
# and are removed, leaving only MAIN. (MAIN's formulas still reference the
# DATA/"Adjusted Expenditure" ranges by name - that is untouched here - so
# switch to manual calculation before deleting the sheets to avoid forcing a
# recalculation that would stamp fresh #REF! values into MAIN's cached
# formula results.)

$wb = $excel.ActiveWorkbook

$excel.Calculation = -4135  # xlCalculationManual
$excel.DisplayAlerts = $false

foreach ($name in @("DATA", "Adjusted Expenditure")) {
    foreach ($sheet in $wb.Worksheets) {
        if ($sheet.Name -eq $name) {
            $sheet.Delete()
        }
    }
}

$excel.DisplayAlerts = $true
